# Commit: "Add files via upload" -- re-upload of RefData_FA_EVENTS.xlsx with
# a couple of small authoring tweaks:
#   1. The 7th sheet tab is renamed from "xxComponentId" to "xxComponentMaster".
#   2. The selection cursor left on that sheet when it was last saved moved
#      from A2 to K27.
# (The workbook-level "first visible tab" scroll hint and the internal
# revisionPtr/document GUID are Excel-session bookkeeping, not something an
# authoring script drives -- they are not reproduced here.)

$wb = $excel.ActiveWorkbook

# 1. Rename the sheet.
$ws = $wb.Worksheets.Item("xxComponentId")
$ws.Name = "xxComponentMaster"

# 2. Make it the active sheet/tab and move the selection to K27, matching
#    the saved cursor position in the diff.
$ws.Activate()
[void]$ws.Range("K27").Select()
